$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "Groups=Smoke" text (shared string used by C2) to "Groups=Retesting"
$ws.Range("C2").Value = "Groups=Retesting"

# Add the new testcase row (row 11): TestCaseNumber 113, Groups "Retesting"
$ws.Range("A11").Value = 113
$ws.Range("A11").HorizontalAlignment = -4131
$ws.Range("B11").Value = "Retesting"
